$d = $word.ActiveDocument

# Locate the paragraph that still contains the placeholder text "Baz" /
# "changes" (split across runs around the _GoBack bookmark).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Baz*chan*") {
        $target = $p
        break
    }
}

$countBefore = $d.Paragraphs.Count

# Build the WordprocessingML fragment for the new paragraph content: a run
# of plain text, then alternating bold "defined terms" and plain connective
# text, all in Arial / dark-grey / white-highlight - matching text copied
# from a web page into Word. A second (initially empty) paragraph is
# appended to carry on the pre-existing "_GoBack" bookmark.
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr      = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="222222"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'
$rPrBold  = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="222222"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'

$body =
  '<w:r>' + $rPr     + '<w:t>A component of software configuration </w:t></w:r>' +
  '<w:r>' + $rPrBold + '<w:t>management</w:t></w:r>' +
  '<w:r>' + $rPr     + '<w:t>,</w:t></w:r>' +
  '<w:r>' + $rPr     + '<w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r>' + $rPrBold + '<w:t>version control</w:t></w:r>' +
  '<w:r>' + $rPr     + '<w:t>, also known as revision </w:t></w:r>' +
  '<w:r>' + $rPrBold + '<w:t>control</w:t></w:r>' +
  '<w:r>' + $rPr     + '<w:t> or source </w:t></w:r>' +
  '<w:r>' + $rPrBold + '<w:t>control</w:t></w:r>' +
  '<w:r>' + $rPr     + '<w:t>, is the </w:t></w:r>' +
  '<w:r>' + $rPrBold + '<w:t>management</w:t></w:r>' +
  '<w:r>' + $rPr     + '<w:t> of changes to documents, computer programs, large web sites, and other collections of information.</w:t></w:r>'

$xml = '<w:p ' + $w + '>' + $body + '</w:p>' +
       '<w:p ' + $w + '><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# Replace the whole old paragraph (text runs + bookmark) with the two new
# paragraphs: the rewritten text, followed by an (empty) paragraph that now
# owns the "_GoBack" bookmark.
[void]$target.Range.InsertXML($xml)

# The insert turned 1 paragraph into 2, so the document now carries one
# extra trailing empty paragraph versus the target; merge the spare empty
# paragraph away. It sits right after the new bookmark paragraph and before
# the document's final (section-ending) paragraph mark, i.e. at the index
# the target paragraph used to occupy (now shifted by the newly inserted
# paragraph).
if ($d.Paragraphs.Count -gt $countBefore) {
    $extraIndex = $countBefore
    $d.Paragraphs.Item($extraIndex).Range.Delete()
}
